$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.208.78"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.882.84"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "479.65"
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.19"
$ws.Range("E6").Value = "  -0.35%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.619"
$ws.Range("E7").Value = "  -2.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.721"
$ws.Range("E9").Value = "  -2.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("E10").Value = "  +7.25%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000351"
$ws.Range("E11").Value = "  +13.81%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.52"
$ws.Range("E12").Value = "  -1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.58"
$ws.Range("E13").Value = "  +1.93%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.520.35"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "14.59"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.902.36"
$ws.Range("E16").Value = "  +1.40%  "
$ws.Range("E17").Value = "  -0.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.66"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.12"
$ws.Range("E19").Value = "  -3.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "68.274.30"
$ws.Range("E20").Value = "  +1.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "434.56"
$ws.Range("E21").Value = "  +0.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "14.67"
$ws.Range("E22").Value = "  -1.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.35"
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "87.72"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.72"
$ws.Range("E25").Value = "  +18.11%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.44"
$ws.Range("E27").Value = "  +3.32%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.96"
$ws.Range("E28").Value = "  +0.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.80"
$ws.Range("E29").Value = "  +4.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "708.94"
$ws.Range("E30").Value = "  -3.34%  "
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.36"
$ws.Range("E31").Value = "  -3.64%  "
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.130"
$ws.Range("E32").Value = "  -2.86%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.85"
$ws.Range("E33").Value = "  +2.23%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0₃0912"
$ws.Range("E34").Value = "  +35.50%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "41.58"
$ws.Range("E35").Value = "  -4.55%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "59.37"
$ws.Range("E36").Value = "  +2.66%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.150"
$ws.Range("E37").Value = "  -6.73%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.64"
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0473"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.08"
$ws.Range("E41").Value = "  +11.28%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.75"
$ws.Range("E42").Value = "  +7.38%  "
$ws.Range("E43").Value = "  +2.91%  "
$ws.Range("E44").Value = "  -1.76%  "
$ws.Range("E45").Value = "  -0.27%  "
$ws.Range("E46").Value = "  -0.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.41"
$ws.Range("E47").Value = "  -1.23%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.13"
$ws.Range("E48").Value = "  -0.92%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "145.98"
$ws.Range("E49").Value = "  +1.18%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "3.13"
$ws.Range("E50").Value = "  -4.66%  "
$ws.Range("E51").Value = "  -2.06%  "
